# Binary Search: Kth smallest price
# Adds row 9 to the "Binary Search 2" sheet: a new problem entry
# ("Kth Smallest Price") with a hyperlink to the Scaler homework page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

# S.no.
$ws.Range("B9").Value = 7

# Question
$ws.Range("D9").Value = "Kth Smallest Price"

# Link -> add hyperlink, then restore the visible cell text
# (TextToDisplay = url keeps the OOXML "display" attribute equal to the
# hyperlink target, matching the other rows in this sheet).
$url = "https://www.scaler.com/academy/mentee-dashboard/class/30365/homework/problems/872?navref=cl_tt_nv"
$h = $ws.Hyperlinks.Add($ws.Range("E9"), $url, [System.Type]::Missing, [System.Type]::Missing, $url)
$ws.Range("E9").Value = "KthPrice - Problem | Scaler Academy"
$ws.Range("E9").WrapText = $true
$ws.Range("E9").HorizontalAlignment = 1
$ws.Range("E9").VerticalAlignment = -4107

# Row height matches the other wrapped rows proportionally (28.8 for a
# single line of wrapped hyperlink text vs 57.6 for two columns wrapping).
$ws.Rows.Item(9).RowHeight = 28.8

# Matches the saved selection/active cell in the edited workbook.
$ws.Range("F9").Select()
